$wb = $excel.ActiveWorkbook

# --- TestSteps sheet ---
$wsTestSteps = $wb.Worksheets.Item("TestSteps")
$wsTestSteps.Activate()

# Fill column H (Results) rows 2 through 19 with "PASS"
for ($r = 2; $r -le 19; $r++) {
    $wsTestSteps.Cells.Item($r, 8).Value = "PASS"
}

# Move the active selection to H20
$wsTestSteps.Range("H20").Select()

# --- TestCases sheet ---
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestCases.Activate()

# Rows 2 and 3: mark "Results" (column C) as "Yes" and add "PASS" in column D
$wsTestCases.Range("C2").Value = "Yes"
$wsTestCases.Range("D2").Value = "PASS"

$wsTestCases.Range("C3").Value = "Yes"
$wsTestCases.Range("D3").Value = "PASS"

# Move the active selection to A2
$wsTestCases.Range("A2").Select()
